# "feat: add 2022-Q1 data"
#
# Before: 2021-Q3, 2021-Q4, 总计(totals)
# After : 2021-Q3, 2021-Q4, 2022-Q1, 总计(totals)
#
#  - the existing "总计" sheet (currently sheet #3) is duplicated; the
#    duplicate keeps the "总计" role (rolled-up totals) and gets a new row
#    for the 2022-Q1 quarter prepended to its table.
#  - the original sheet #3 is repurposed to hold the detailed per-fund
#    holdings for the new "2022-Q1" quarter (same layout as the
#    "2021-Q3"/"2021-Q4" sheets) and is renamed "2022-Q1".

# Helper: write $Text into $CellRef as a genuine text value (not a
# number), without leaving a lingering "quote prefix" style behind.
# We do this by entering the value with a leading apostrophe (which
# forces Excel to treat it as text even when it looks numeric) and then
# re-applying the cell's original formatting from a known-blank cell so
# no extra style/quote-prefix survives.
function Set-TextValue {
    param($Sheet, $CellRef, $Text, $Scratch)
    $Scratch.Copy()
    $target = $Sheet.Range($CellRef)
    $target.Value = "'" + $Text
    $target.PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(3)   # currently named "总计"

# 1) Duplicate the "总计" sheet; the copy is placed right after it and
#    will become the new, updated "总计" sheet (sheetId 4). Rename the
#    original out of the way first so the two sheet names don't collide.
$totalSheet.Name = "2022-Q1-tmp"
$totalSheet.Copy($null, $totalSheet)
$newTotalSheet = $wb.Worksheets.Item(4)
$newTotalSheet.Name = "总计"

# 2) Insert a new first data row (2022-Q1) into the duplicated "总计"
#    sheet, pushing the existing rows (2021-Q4, 2021-Q3) down by one.
#    Row 4 is brand-new, so give it the same style as the existing
#    index-column cells (A2/A3) before writing into it.
$newTotalSheet.Range("A2").Copy()
$newTotalSheet.Range("A4").PasteSpecial(-4122)

$newTotalSheet.Range("B4").Value = "2021-Q3"
$newTotalSheet.Range("C4").Value = 2
$newTotalSheet.Range("D4").Value = 0.3
$newTotalSheet.Range("A4").Value = 2

$newTotalSheet.Range("B3").Value = "2021-Q4"
$newTotalSheet.Range("C3").Value = 2
$newTotalSheet.Range("D3").Value = 0.28
$newTotalSheet.Range("A3").Value = 1

$newTotalSheet.Range("B2").Value = "2022-Q1"
$newTotalSheet.Range("C2").Value = 2
$newTotalSheet.Range("D2").Value = 0.33
$newTotalSheet.Range("A2").Value = 0

# 3) Repurpose the original "总计" sheet (position 3) into the detailed
#    "2022-Q1" fund-holdings sheet, matching the layout used by the
#    "2021-Q3" / "2021-Q4" sheets.
$q1Sheet = $totalSheet
$scratch = $q1Sheet.Range("Z100")

# Copy the header cell formatting (style already used by this sheet)
# across the new columns E:H so they pick up the same style as B1:D1.
$q1Sheet.Range("B1").Copy()
$q1Sheet.Range("E1:H1").PasteSpecial(-4122)

$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

Set-TextValue $q1Sheet "B2" "513030" $scratch
$q1Sheet.Range("C2").Value = "华安国际龙头(DAX)ETFQDII"
Set-TextValue $q1Sheet "D2" "6.49" $scratch
Set-TextValue $q1Sheet "E2" "92.80" $scratch
Set-TextValue $q1Sheet "F2" "4.69" $scratch
Set-TextValue $q1Sheet "G2" "0.3044" $scratch
$q1Sheet.Range("H2").Value = 5

Set-TextValue $q1Sheet "B3" "513080" $scratch
$q1Sheet.Range("C3").Value = "华安法国CAC40ETF（QDII）"
Set-TextValue $q1Sheet "D3" "0.60" $scratch
Set-TextValue $q1Sheet "E3" "96.69" $scratch
Set-TextValue $q1Sheet "F3" "3.98" $scratch
Set-TextValue $q1Sheet "G3" "0.0239" $scratch
$q1Sheet.Range("H3").Value = 7

$scratch.Clear()

$q1Sheet.Name = "2022-Q1"
